# Update "想去人数" (interest count) values in the "展览" and "全部类型" sheets
# to reflect the newly scraped figures (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 1450
$ws1.Range("F5").Value  = 12100
$ws1.Range("F6").Value  = 4471
$ws1.Range("F8").Value  = 60
$ws1.Range("F12").Value = 1118
$ws1.Range("F14").Value = 61
$ws1.Range("F15").Value = 5239
$ws1.Range("F20").Value = 11479

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 1450
$ws4.Range("F5").Value  = 12100
$ws4.Range("F6").Value  = 4471
$ws4.Range("F8").Value  = 60
$ws4.Range("F13").Value = 1118
$ws4.Range("F15").Value = 61
$ws4.Range("F16").Value = 5239
$ws4.Range("F21").Value = 11479

$wb.Save()
